# Edit script: implements the rdbms assertResultMatch/assertResultNotMatch
# commands, the step.inTime (waitMs variants of step) commands, the
# localdb queryAsCSV command, and the web assertElementDisabled /
# checkAll(locator,waitMs) / uncheckAll(locator,waitMs) commands, by
# updating the hidden "#system" command-reference sheet and the
# associated named ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- rdbms (column T): add assertResultMatch / assertResultNotMatch ---
$tList = @(
    "assertResultMatch(var,columns,search)",
    "assertResultNotMatch(var,columns,search)",
    "resultToCSV(var,csvFile,delim,showHeader)",
    "runFile(var,db,file)",
    "runSQL(var,db,sql)",
    "runSQLs(var,db,sqls)",
    "saveResult(db,sql,output)",
    "saveResults(db,sqls,outputDir)"
)
$row = 2
foreach ($item in $tList) {
    $ws.Cells.Item($row, 20).Value = $item
    $row = $row + 1
}

# --- localdb (column O): add queryAsCSV(var,sql) ---
$oList = @(
    "cloneTable(var,source,target)",
    "dropTables(var,tables)",
    "exportCSV(sql,output)",
    "exportEXCEL(sql,output,sheet,start)",
    "exportJSON(sql,output,header)",
    "exportXML(sql,output,root,row,cell)",
    "importCSV(var,csv,table)",
    "importEXCEL(var,excel,sheet,ranges,table)",
    "importRecords(var,sourceDb,sql,table)",
    "purge(var)",
    "queryAsCSV(var,sql)",
    "runSQLs(var,sqls)"
)
$row = 2
foreach ($item in $oList) {
    $ws.Cells.Item($row, 15).Value = $item
    $row = $row + 1
}

# --- step.inTime (column Z): replace old tn.5250 data with the new
# step.inTime header + observe/perform/validate (waitMs variants) ---
$ws.Cells.Item(1, 26).Value = "step.inTime"
$zList = @(
    "observe(prompt,waitMs)",
    "perform(instructions,waitMs)",
    "validate(prompt,responses,passResponses,waitMs)"
)
$row = 2
foreach ($item in $zList) {
    $ws.Cells.Item($row, 26).Value = $item
    $row = $row + 1
}
$ws.Cells.Item(5, 26).ClearContents()
$ws.Cells.Item(6, 26).ClearContents()

# --- target (column A): rename the "tn.5250" entry to "step.inTime" ---
$ws.Cells.Item(26, 1).Value = "step.inTime"

# --- web (column AA): insert assertElementDisabled(locator),
# replace checkAll(locator)/uncheckAll(locator) with the new
# waitMs-enabled signatures, and re-sort the full command list ---
$aaList = @(
    "assertAndClick(locator,label)",
    "assertAttribute(locator,attrName,value)",
    "assertAttributeContain(locator,attrName,contains)",
    "assertAttributeNotContain(locator,attrName,contains)",
    "assertAttributeNotPresent(locator,attrName)",
    "assertAttributePresent(locator,attrName)",
    "assertChecked(locator)",
    "assertContainCount(locator,text,count)",
    "assertCssNotPresent(locator,property)",
    "assertCssPresent(locator,property,value)",
    "assertElementByAttributes(nameValues)",
    "assertElementByText(locator,text)",
    "assertElementCount(locator,count)",
    "assertElementDisabled(locator)",
    "assertElementEnabled(locator)",
    "assertElementNotPresent(locator)",
    "assertElementPresent(locator)",
    "assertElementsPresent(prefix)",
    "assertFocus(locator)",
    "assertFrameCount(count)",
    "assertFramePresent(frameName)",
    "assertIECompatMode()",
    "assertIENativeMode()",
    "assertLinkByLabel(label)",
    "assertMultiSelect(locator)",
    "assertNotChecked(locator)",
    "assertNotFocus(locator)",
    "assertNotText(locator,text)",
    "assertNotVisible(locator)",
    "assertOneMatch(locator)",
    "assertScrollbarHNotPresent(locator)",
    "assertScrollbarHPresent(locator)",
    "assertScrollbarVNotPresent(locator)",
    "assertScrollbarVPresent(locator)",
    "assertSingleSelect(locator)",
    "assertTable(locator,row,column,text)",
    "assertText(locator,text)",
    "assertTextContains(locator,text)",
    "assertTextCount(locator,text,count)",
    "assertTextList(locator,list,ignoreOrder)",
    "assertTextMatches(text,minMatch,scrollTo)",
    "assertTextNotContain(locator,text)",
    "assertTextNotPresent(text)",
    "assertTextOrder(locator,descending)",
    "assertTextPresent(text)",
    "assertTitle(text)",
    "assertValue(locator,value)",
    "assertValueOrder(locator,descending)",
    "assertVisible(locator)",
    "checkAll(locator,waitMs)",
    "clearLocalStorage()",
    "click(locator)",
    "clickAll(locator)",
    "clickAndWait(locator,waitMs)",
    "clickByLabel(label)",
    "clickByLabelAndWait(label,waitMs)",
    "clickOffset(locator,x,y)",
    "clickWithKeys(locator,keys)",
    "close()",
    "closeAll()",
    "deselect(locator,text)",
    "deselectMulti(locator,array)",
    "dismissInvalidCert()",
    "dismissInvalidCertPopup()",
    "doubleClick(locator)",
    "doubleClickAndWait(locator,waitMs)",
    "doubleClickByLabel(label)",
    "doubleClickByLabelAndWait(label,waitMs)",
    "dragAndDrop(fromLocator,toLocator)",
    "dragTo(fromLocator,xOffset,yOffset)",
    "editLocalStorage(key,value)",
    "executeScript(var,script)",
    "focus(locator)",
    "goBack()",
    "goBackAndWait()",
    "maximizeWindow()",
    "mouseOver(locator)",
    "open(url)",
    "openAndWait(url,waitMs)",
    "openHttpBasic(url,username,password)",
    "openIgnoreTimeout(url)",
    "refresh()",
    "refreshAndWait()",
    "resizeWindow(width,height)",
    "rightClick(locator)",
    "saveAllWindowIds(var)",
    "saveAllWindowNames(var)",
    "saveAttribute(var,locator,attrName)",
    "saveAttributeList(var,locator,attrName)",
    "saveBrowserVersion(var)",
    "saveCount(var,locator)",
    "saveDivsAsCsv(headers,rows,cells,nextPage,file)",
    "saveElement(var,locator)",
    "saveElements(var,locator)",
    "saveInfiniteDivsAsCsv(config,file)",
    "saveInfiniteTableAsCsv(config,file)",
    "saveLocalStorage(var,key)",
    "saveLocation(var)",
    "savePageAs(var,sessionIdName,url)",
    "savePageAsFile(sessionIdName,url,file)",
    "saveSelectedText(var,locator)",
    "saveSelectedValue(var,locator)",
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveTitle(var)",
    "saveValue(var,locator)",
    "saveValues(var,locator)",
    "screenshot(file,locator)",
    "scrollElement(locator,xOffset,yOffset)",
    "scrollLeft(locator,pixel)",
    "scrollPage(xOffset,yOffset)",
    "scrollRight(locator,pixel)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectAllOptions(locator)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiByValue(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "switchBrowser(profile,config)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator,waitMs)",
    "unselectAllText()",
    "updateAttribute(locator,attrName,value)",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator)",
    "waitForElementsPresent(locators)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)"
)
$row = 2
foreach ($item in $aaList) {
    $ws.Cells.Item($row, 27).Value = $item
    $row = $row + 1
}

# --- defined names: widen localdb/rdbms/web ranges, add step.inTime ---
$wb.Names.Item("localdb").RefersTo = "='#system'!`$O`$2:`$O`$13"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$145"
$wb.Names.Add("step.inTime", "='#system'!`$Z`$2:`$Z`$4")

